$d = $word.ActiveDocument

# The bibliography paragraph currently has the 5 numbered references
# concatenated together in a single run/text. Split them onto separate
# lines using line breaks (w:br), using Find/Replace with wildcards so the
# line break is inserted exactly at each "<year/punct>N)" boundary without
# touching anything else in the document.

$d.Content.Find.Execute("2005.2)", $false, $false, $true, $false, $false, $true, 1, $false, "2005.^l2)", 2) | Out-Null
$d.Content.Find.Execute("2018.3)", $false, $false, $true, $false, $false, $true, 1, $false, "2018.^l3)", 2) | Out-Null
$d.Content.Find.Execute("2012.4)", $false, $false, $true, $false, $false, $true, 1, $false, "2012.^l4)", 2) | Out-Null
$d.Content.Find.Execute("2014.5)", $false, $false, $true, $false, $false, $true, 1, $false, "2014.^l5)", 2) | Out-Null
